# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Update the K column (G) values for each outing row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value (column G)
$kValues = @{
    2  = 0
    3  = 2
    5  = 2
    6  = 0
    7  = 1
    8  = 2
    9  = 2
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 2
    17 = 0
    18 = 3
    19 = 0
    20 = 1
    21 = 1
    22 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
